# The MATLAB "save results" routine now appends each Monte Carlo run's
# results as a new row instead of overwriting row 2 each time. This
# simulates three additional runs being appended to the Results sheet
# (rows 3, 4 and 5), each with the same input parameters as row 2 but
# fresh Lead Time / Idle Time statistics, leaving the summary columns
# (U:X) blank for the appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3, @(4.9, 10, 11, 15, 16, 1, 1, 1, 1, 0.04, 0.04, 0.04, 0.04, 200, 300, 600, 0.418, 0.33, 0.5, 50)),
    @(4, @(4.9, 10, 11, 15, 16, 1, 1, 1, 1, 0.04, 0.04, 0.04, 0.04, 200, 300, 600, 0.418, 0.33, 0.5, 50)),
    @(5, @(5.1, 10, 11, 15, 16, 1, 1, 1, 1, 0.04, 0.04, 0.04, 0.04, 200, 300, 600, 0.175, 0.33, 0.5, 50))
)

foreach ($entry in $newRows) {
    $rowNum = $entry[0]
    $values = $entry[1]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($rowNum, $col + 1).Value = $values[$col]
    }
}

# Excel leaves the newly written range selected after the save.
$ws.Range("A5:T5").Select()
